$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("high_loadings")
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(5, 2).Value = 2
$ws.Cells.Item(6, 2).Value = 2
$ws.Cells.Item(8, 2).Value = 3
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(11, 2).Value = 3
$ws.Cells.Item(13, 2).Value = 3
$ws.Cells.Item(14, 2).Value = 3
$ws.Cells.Item(15, 2).Value = 3
$ws.Cells.Item(16, 2).Value = 2
$ws.Cells.Item(17, 2).Value = 2
$ws.Cells.Item(18, 2).Value = 1
$ws.Cells.Item(20, 2).Value = 2
$ws.Cells.Item(21, 2).Value = 2
$ws.Cells.Item(22, 2).Value = 3
$ws.Cells.Item(25, 2).Value = 3
$ws.Cells.Item(26, 2).Value = 2
$ws.Cells.Item(28, 2).Value = 2
$ws.Cells.Item(29, 2).Value = 2
$ws.Cells.Item(30, 2).Value = 2
$ws.Cells.Item(31, 2).Value = 3
$ws.Cells.Item(33, 2).Value = 3
$ws.Cells.Item(34, 2).Value = 2
$ws.Cells.Item(36, 2).Value = 3
$ws.Cells.Item(37, 2).Value = 2
$ws.Cells.Item(38, 2).Value = 1
$ws.Cells.Item(39, 2).Value = 1
$ws.Cells.Item(40, 2).Value = 1
$ws.Cells.Item(41, 2).Value = 3
$ws.Cells.Item(43, 2).Value = 2
$ws.Cells.Item(44, 2).Value = 3
$ws.Cells.Item(45, 2).Value = 1
$ws.Cells.Item(46, 2).Value = 1
$ws.Cells.Item(47, 2).Value = 3
$ws.Cells.Item(48, 2).Value = 3
$ws.Cells.Item(49, 2).Value = 1
$ws.Cells.Item(50, 2).Value = 2
$ws.Cells.Item(51, 2).Value = 2
$ws.Cells.Item(52, 2).Value = 2
$ws.Cells.Item(53, 2).Value = 1
$ws.Cells.Item(54, 2).Value = 1
$ws.Cells.Item(55, 2).Value = 1
$ws.Cells.Item(56, 2).Value = 1
$ws.Cells.Item(58, 2).Value = 2
$ws.Cells.Item(60, 2).Value = 2
$ws.Cells.Item(61, 2).Value = 1
$ws.Cells.Item(62, 2).Value = 1
$ws.Cells.Item(63, 2).Value = 2
$ws.Cells.Item(66, 2).Value = 3
$ws.Cells.Item(67, 2).Value = 1
$ws.Cells.Item(69, 2).Value = 3

$ws = $wb.Worksheets.Item("Influencers_uniques")
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(5, 2).Value = 2
$ws.Cells.Item(6, 2).Value = 2
$ws.Cells.Item(7, 2).Value = 3
$ws.Cells.Item(9, 2).Value = 3
$ws.Cells.Item(12, 2).Value = 2
$ws.Cells.Item(13, 2).Value = 2
$ws.Cells.Item(14, 2).Value = 3
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(18, 2).Value = 3
$ws.Cells.Item(21, 2).Value = 2
$ws.Cells.Item(24, 2).Value = 1
$ws.Cells.Item(26, 2).Value = 2
$ws.Cells.Item(27, 2).Value = 1
$ws.Cells.Item(28, 2).Value = 3
$ws.Cells.Item(29, 2).Value = 3
$ws.Cells.Item(30, 2).Value = 3
$ws.Cells.Item(31, 2).Value = 1
$ws.Cells.Item(32, 2).Value = 2
$ws.Cells.Item(33, 2).Value = 2
$ws.Cells.Item(35, 2).Value = 1
$ws.Cells.Item(36, 2).Value = 3
$ws.Cells.Item(38, 2).Value = 1
$ws.Cells.Item(39, 2).Value = 2
$ws.Cells.Item(40, 2).Value = 3
$ws.Cells.Item(42, 2).Value = 1
$ws.Cells.Item(43, 2).Value = 1
$ws.Cells.Item(46, 2).Value = 3
$ws.Cells.Item(47, 2).Value = 1
$ws.Cells.Item(48, 2).Value = 1
$ws.Cells.Item(50, 2).Value = 1
$ws.Cells.Item(51, 2).Value = 3
$ws.Cells.Item(53, 2).Value = 3
$ws.Cells.Item(54, 2).Value = 3
$ws.Cells.Item(57, 2).Value = 2
$ws.Cells.Item(58, 2).Value = 1
$ws.Cells.Item(60, 2).Value = 1
$ws.Cells.Item(62, 2).Value = 2
$ws.Cells.Item(63, 2).Value = 2
$ws.Cells.Item(65, 2).Value = 3
$ws.Cells.Item(67, 2).Value = 3
$ws.Cells.Item(70, 2).Value = 2
$ws.Cells.Item(71, 2).Value = 3
